$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1017.1964
$ws.Range("I15").Value = 1017.1964
$ws.Range("K15").Value = 3051.5892
$ws.Range("M15").Value = -2882.5892

$ws.Range("H43").Value = 1001
$ws.Range("J43").Value = 1001
$ws.Range("L43").Value = 1001
$ws.Range("N43").Value = -1139

$ws.Range("H70").Value = 4786.079
$ws.Range("I70").Value = 2147.6365
$ws.Range("J70").Value = 22199.8
$ws.Range("K70").Value = 6442.9095
$ws.Range("L70").Value = 66599.39999999999
$ws.Range("M70").Value = -6172.9095
$ws.Range("N70").Value = -67139.39999999999

$ws.Range("H73").Value = 4786.079
$ws.Range("I73").Value = 2147.6365
$ws.Range("J73").Value = 22199.8
$ws.Range("K73").Value = 6442.9095
$ws.Range("L73").Value = 66599.39999999999
$ws.Range("M73").Value = -5506.9095
$ws.Range("N73").Value = -68471.39999999999

$ws.Range("H92").Value = 12500256
$ws.Range("I92").Value = 13889063
$ws.Range("K92").Value = 13889063
$ws.Range("M92").Value = -13887815

$ws.Range("H96").Value = 1555.6
$ws.Range("I96").Value = 1469
$ws.Range("J96").Value = 1685.5
$ws.Range("K96").Value = 4407
$ws.Range("L96").Value = 5056.5
$ws.Range("M96").Value = -3034
$ws.Range("N96").Value = -7802.5

$ws.Range("H100").Value = 1652.8889
$ws.Range("I100").Value = 900.25
$ws.Range("K100").Value = 900.25
$ws.Range("M100").Value = -359.25

$ws.Range("H125").Value = 365.5
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 365.5
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 3289.5
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -8209.5

$ws.Range("H132").Value = 1237.1111
$ws.Range("I132").Value = 1142.1904
$ws.Range("K132").Value = 3426.5712
$ws.Range("M132").Value = -896.5711999999999

$ws.Range("I135").Value = 1002.1429
$ws.Range("K135").Value = 9019.286100000001
$ws.Range("M135").Value = -6484.286100000001

$ws.Range("H136").Value = 66513.8
$ws.Range("J136").Value = 66513.8
$ws.Range("L136").Value = 66513.8
$ws.Range("N136").Value = -76713.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4979.939
$ws.Range("I32").Value = 3964.8333
$ws.Range("K32").Value = 3964.8333
$ws.Range("M32").Value = -3677.8333

$ws.Range("H51").Value = 30000
$ws.Range("J51").Value = 30000
$ws.Range("L51").Value = 30000
$ws.Range("N51").Value = -31512

$ws.Range("H97").Value = 1095.5294
$ws.Range("I97").Value = 1066.8125
$ws.Range("K97").Value = 1066.8125
$ws.Range("M97").Value = -570.8125

$ws.Range("H102").Value = 1785.8
$ws.Range("I102").Value = 1785.8
$ws.Range("K102").Value = 1785.8
$ws.Range("M102").Value = -163.8

$ws.Range("H131").Value = 48998.363
$ws.Range("J131").Value = 48998.363
$ws.Range("L131").Value = 48998.363
$ws.Range("N131").Value = -59078.363

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1228.8572
$ws.Range("I99").Value = 874
$ws.Range("K99").Value = 874
$ws.Range("M99").Value = 624

$ws.Range("H134").Value = 5878.5
$ws.Range("I134").Value = 6350.696
$ws.Range("J134").Value = 2258.3333
$ws.Range("K134").Value = 19052.088
$ws.Range("L134").Value = 6774.999899999999
$ws.Range("M134").Value = -16517.088
$ws.Range("N134").Value = -11844.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2741.9473
$ws.Range("I31").Value = 2636.2222
$ws.Range("J31").Value = 2837.1
$ws.Range("K31").Value = 2636.2222
$ws.Range("L31").Value = 2837.1
$ws.Range("M31").Value = -2341.2222
$ws.Range("N31").Value = -3427.1

$ws.Range("H34").Value = 2741.9473
$ws.Range("I34").Value = 2636.2222
$ws.Range("J34").Value = 2837.1
$ws.Range("K34").Value = 2636.2222
$ws.Range("L34").Value = 2837.1
$ws.Range("M34").Value = -2434.2222
$ws.Range("N34").Value = -3241.1

$ws.Range("H62").Value = 2103.5386
$ws.Range("I62").Value = 2219.5
$ws.Range("J62").Value = 1918
$ws.Range("K62").Value = 2219.5
$ws.Range("L62").Value = 1918
$ws.Range("M62").Value = -1595.5
$ws.Range("N62").Value = -3166

$ws.Range("H65").Value = 2103.5386
$ws.Range("I65").Value = 2219.5
$ws.Range("J65").Value = 1918
$ws.Range("K65").Value = 11097.5
$ws.Range("L65").Value = 9590
$ws.Range("M65").Value = -7977.5
$ws.Range("N65").Value = -15830

$ws.Range("H122").Value = 5226.8184
$ws.Range("I122").Value = 3640.5715
$ws.Range("J122").Value = 8002.75
$ws.Range("K122").Value = 10921.7145
$ws.Range("L122").Value = 24008.25
$ws.Range("M122").Value = -8471.7145
$ws.Range("N122").Value = -28908.25

$ws.Range("H132").Value = 3219
$ws.Range("I132").Value = 2507.8667
$ws.Range("K132").Value = 7523.6001
$ws.Range("M132").Value = -4993.6001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 500
$ws.Range("J97").Value = 500
$ws.Range("L97").Value = 1500
$ws.Range("N97").Value = -2492

$ws.Range("H131").Value = 13923.227
$ws.Range("I131").Value = 696
$ws.Range("J131").Value = 15301.0625
$ws.Range("K131").Value = 2088
$ws.Range("L131").Value = 45903.1875
$ws.Range("M131").Value = 2952
$ws.Range("N131").Value = -55983.1875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 28803
$ws.Range("I13").Value = 600
$ws.Range("K13").Value = 600
$ws.Range("M13").Value = -461

$ws.Range("H70").Value = 3666.5557
$ws.Range("I70").Value = 3333.3333
$ws.Range("K70").Value = 3333.3333
$ws.Range("M70").Value = -3063.3333

$ws.Range("H73").Value = 3666.5557
$ws.Range("I73").Value = 3333.3333
$ws.Range("K73").Value = 3333.3333
$ws.Range("M73").Value = -2397.3333

$ws.Range("H99").Value = 8599
$ws.Range("I99").Value = 8599
$ws.Range("K99").Value = 8599
$ws.Range("M99").Value = -6353

$ws.Range("H102").Value = 3470.3872
$ws.Range("I102").Value = 3885.45
$ws.Range("K102").Value = 3885.45
$ws.Range("M102").Value = -2263.45

$ws.Range("H113").Value = 1262.6364
$ws.Range("I113").Value = 938.5
$ws.Range("J113").Value = 1447.8572
$ws.Range("K113").Value = 938.5
$ws.Range("L113").Value = 1447.8572
$ws.Range("M113").Value = 1231.5
$ws.Range("N113").Value = -5787.8572

$ws.Range("H122").Value = 1068.3214
$ws.Range("I122").Value = 1070.8889
$ws.Range("J122").Value = 999
$ws.Range("K122").Value = 3212.6667
$ws.Range("L122").Value = 2997
$ws.Range("M122").Value = -762.6666999999998
$ws.Range("N122").Value = -7897

$ws.Range("H132").Value = 1849.2894
$ws.Range("I132").Value = 1360.2307
$ws.Range("K132").Value = 4080.6921
$ws.Range("M132").Value = -1550.6921

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5054.5557
$ws.Range("I7").Value = 3920.8
$ws.Range("K7").Value = 3920.8
$ws.Range("M7").Value = -3808.8

$ws.Range("H40").Value = 19787.334
$ws.Range("I40").Value = 30673
$ws.Range("K40").Value = 30673
$ws.Range("M40").Value = -30537

$ws.Range("H100").Value = 1624.75
$ws.Range("I100").Value = 1499.6666
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 1499.6666
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -958.6666
$ws.Range("N100").Value = -3082

$ws.Range("H122").Value = 5089.0713
$ws.Range("I122").Value = 4905.875
$ws.Range("K122").Value = 14717.625
$ws.Range("M122").Value = -12267.625

$ws.Range("H126").Value = 5054.5557
$ws.Range("I126").Value = 3920.8
$ws.Range("K126").Value = 11762.4
$ws.Range("M126").Value = -9292.400000000001

$ws.Range("H132").Value = 5834.1333
$ws.Range("I132").Value = 6232.8887
$ws.Range("J132").Value = 5236
$ws.Range("K132").Value = 18698.6661
$ws.Range("L132").Value = 15708
$ws.Range("M132").Value = -16168.6661
$ws.Range("N132").Value = -20768

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()

$ws.Range("H122").Value = 79476.60000000001
$ws.Range("I122").Value = 131177.83
$ws.Range("J122").Value = 1924.75
$ws.Range("K122").Value = 393533.49
$ws.Range("L122").Value = 5774.25
$ws.Range("M122").Value = -391083.49
$ws.Range("N122").Value = -10674.25

$ws.Range("H132").Value = 2324.842
$ws.Range("I132").Value = 1368.3334
$ws.Range("J132").Value = 3964.5715
$ws.Range("K132").Value = 4105.0002
$ws.Range("L132").Value = 11893.7145
$ws.Range("M132").Value = -1575.0002
$ws.Range("N132").Value = -16953.7145
